{"js": "// Ajusta ponto m\u00e9dio quest\u00e3o 1\n// The \"Xm\" (ponto m\u00e9dio / midpoint) column of the frequency table held the\n// wrong values (copy/paste of the \"fi acumulada\" column). Replace each\n// midpoint cell with the correct computed class midpoint.\nconst replacements = [\n  [\"225\", \"14,16\"],\n  [\"227\", \"14,48\"],\n  [\"229\", \"14,80\"],\n  [\"231\", \"15,12\"],\n  [\"233\", \"15,44\"],\n  [\"235\", \"15,76\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldVal, newVal] of replacements) {\n  const results = body.search(oldVal, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly one match for '\" + oldVal + \"', found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(newVal, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Ajusta ponto m\u00e9dio quest\u00e3o 1\n# The \"Xm\" (ponto m\u00e9dio / midpoint) column of the frequency table held the\n# wrong values (copy/paste of the \"fi acumulada\" column). Replace each\n# midpoint cell with the correct computed class midpoint.\n\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"225\" = \"14,16\"\n    \"227\" = \"14,48\"\n    \"229\" = \"14,80\"\n    \"231\" = \"15,12\"\n    \"233\" = \"15,44\"\n    \"235\" = \"15,76\"\n}\n\nforeach ($oldVal in $replacements.Keys) {\n    $newVal = $replacements[$oldVal]\n    $matched = $false\n\n    foreach ($p in $d.Paragraphs) {\n        $rng = $p.Range\n        $text = $rng.Text\n\n        # Paragraph text includes the trailing paragraph mark (chr 13) and,\n        # inside a table cell, the cell mark (chr 7); strip those before\n        # comparing so we match a cell whose whole content is just the value.\n        $trimmed = $text -replace \"[\\r\\a]+$\", \"\"\n\n        if ($trimmed -eq $oldVal) {\n            [void]$rng.MoveEnd(1, -($text.Length - $trimmed.Length))\n            $rng.Text = $newVal\n            $matched = $true\n            break\n        }\n    }\n\n    if (-not $matched) {\n        throw \"Could not find a unique paragraph/cell containing '$oldVal'\"\n    }\n}\n"}
